$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 7665
$ws1.Range("F6").Value = 5577
$ws1.Range("F7").Value = 458
$ws1.Range("F8").Value = 71
$ws1.Range("F11").Value = 249
$ws1.Range("F12").Value = 203
$ws1.Range("F13").Value = 51

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 7665
$ws4.Range("F6").Value = 5577
$ws4.Range("F7").Value = 458
$ws4.Range("F8").Value = 71
$ws4.Range("F11").Value = 249
$ws4.Range("F14").Value = 203
$ws4.Range("F15").Value = 51
